$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '35.426.47'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.80%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.919.56'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +3.40%  '
$ws.Range('E4').Value = '  -0.50%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.45'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.25%  '
$ws.Range('E6').Value = '  +6.22%  '
$ws.Range('E7').Value = '  -0.50%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.65'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.41%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.352'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +7.23%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '52.66'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +12.23%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0718'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.63%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0996'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.61%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.190.63'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.03%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '12.11'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.65%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.700'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.60%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.914.66'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.17%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.90'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.90%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '35.443.23'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.88%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '72.10'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.42%  '
$ws.Range('E20').Value = '  +3.64%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '240.59'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.53'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.55%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.87'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.03%  '
$ws.Range('E24').Value = '  -0.43%  '
$ws.Range('E25').Value = '  +1.47%  '
$ws.Range('E26').Value = '  +21.49%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '170.60'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.47'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.21%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.53'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.95%  '
$ws.Range('E31').Value = '  +4.02%  '
$ws.Range('E32').Value = '  +1.87%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.943'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +14.16%  '
$ws.Range('B34').Value = 'BinanceUSD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.01'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.39%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.14'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.29%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.74'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.68%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.04'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.91%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.33'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.93%  '
$ws.Range('E39').Value = '  +2.14%  '
$ws.Range('E40').Value = '  +4.18%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0653'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +17.13%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '16.28'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +9.09%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '90.52'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.96%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.345.14'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.24%  '
$ws.Range('E45').Value = '  +2.74%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '48.38'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +38.98%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.81'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.34%  '
$ws.Range('E48').Value = '  -0.10%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.59'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.25%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.102.79'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.16%  '
$ws.Range('E51').Value = '  +2.88%  '
